$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,6
$row2[0,0] = 0.1783908196033299
$row2[0,1] = 0.3606156554386025
$row2[0,2] = 0.2599511937740667
$row2[0,3] = 0.5098540906711122
$row2[0,4] = 0.4943913024279584
$row2[0,5] = 15
$ws.Range("B2:G2").Value = $row2

$row3 = New-Object 'object[,]' 1,6
$row3[0,0] = 0.2147957804815247
$row3[0,1] = 0.3580819848962541
$row3[0,2] = 0.2088379301833887
$row3[0,3] = 0.4569878884427777
$row3[0,4] = 0.4185882721160311
$row3[0,5] = 14
$ws.Range("B3:G3").Value = $row3

$row4 = New-Object 'object[,]' 1,6
$row4[0,0] = 0.245473867636748
$row4[0,1] = 0.3683843928281823
$row4[0,2] = 0.2358514199634523
$row4[0,3] = 0.4856453644002507
$row4[0,4] = 0.4361500127557235
$row4[0,5] = 13
$ws.Range("B4:G4").Value = $row4

$row5 = New-Object 'object[,]' 1,6
$row5[0,0] = 0.301463903713677
$row5[0,1] = 0.3567050099425655
$row5[0,2] = 0.2583510523008795
$row5[0,3] = 0.5082824532687308
$row5[0,4] = 0.4274285484895888
$row5[0,5] = 12
$ws.Range("B5:G5").Value = $row5

$row6 = New-Object 'object[,]' 1,6
$row6[0,0] = 0.2900334324511264
$row6[0,1] = 0.3439303588265401
$row6[0,2] = 0.1874043702602189
$row6[0,3] = 0.4329022640968963
$row6[0,4] = 0.3370659818980856
$row6[0,5] = 11
$ws.Range("B6:G6").Value = $row6

$row7 = New-Object 'object[,]' 1,6
$row7[0,0] = 0.2833661826906564
$row7[0,1] = 0.2909807410604855
$row7[0,2] = 0.2099604875429734
$row7[0,3] = 0.4582144558424291
$row7[0,4] = 0.3795671424286398
$row7[0,5] = 10
$ws.Range("B7:G7").Value = $row7

$row8 = New-Object 'object[,]' 1,6
$row8[0,0] = 0.3070978867771534
$row8[0,1] = 0.3733636115588901
$row8[0,2] = 0.2272551859511885
$row8[0,3] = 0.476712896774556
$row8[0,4] = 0.3867354821117653
$row8[0,5] = 9
$ws.Range("B8:G8").Value = $row8

$row9 = New-Object 'object[,]' 1,6
$row9[0,0] = 0.3125739411747067
$row9[0,1] = 0.4101416856701894
$row9[0,2] = 0.3161064439746084
$row9[0,3] = 0.5622334425971194
$row9[0,4] = 0.499604386659364
$row9[0,5] = 8
$ws.Range("B9:G9").Value = $row9

$row10 = New-Object 'object[,]' 1,6
$row10[0,0] = 0.2627214152773011
$row10[0,1] = 0.3235434349789264
$row10[0,2] = 0.1669924024756296
$row10[0,3] = 0.4086470389904099
$row10[0,4] = 0.3380801243818421
$row10[0,5] = 7
$ws.Range("B10:G10").Value = $row10

$row11 = New-Object 'object[,]' 1,6
$row11[0,0] = 0.3023013051968123
$row11[0,1] = 0.3856149542229345
$row11[0,2] = 0.3827274330185637
$row11[0,3] = 0.618649685216572
$row11[0,4] = 0.5912779588939882
$row11[0,5] = 6
$ws.Range("B11:G11").Value = $row11

